# The deck ships with two theme parts:
#   ppt/theme/theme1.xml -> "Integral" design / "Red Violet" colour scheme,
#                            wired to the slide master (i.e. every slide).
#   ppt/theme/theme2.xml -> "Office Theme" design / "Office" colour scheme,
#                            wired to the notes master.
#
# The authored change swaps the two themes' contents: the slide master's
# theme becomes the "Office Theme" colour scheme (dk1/lt1/dk2/lt2/accent1-6/
# hlink/folHlink), and the notes master's theme becomes the former
# "Integral"/"Red Violet" colour scheme.
#
# Apply the reachable half of that swap through the object model: push the
# "Office Theme" palette onto the presentation's (slide master's) theme via
# Theme.ThemeColorScheme - the supported, documented way to edit a theme's
# colours from PowerPoint automation.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

# ThemeColorScheme.Item(n).RGB uses the standard OLE colour encoding
# (R + G*256 + B*65536). Index order: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1-6, 11 hlink, 12 folHlink.
$officeThemeColors = @(
    0,        # dk1      000000
    16777215, # lt1      FFFFFF
    6968388,  # dk2      44546A
    15132391, # lt2      E7E6E6
    13998939, # accent1  5B9BD5
    3243501,  # accent2  ED7D31
    10855845, # accent3  A5A5A5
    49407,    # accent4  FFC000
    12874308, # accent5  4472C4
    4697456,  # accent6  70AD47
    12673797, # hlink    0563C1
    7491477   # folHlink 954F72
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeThemeColors[$i - 1]
}
